$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters -> new values for rows 2 and 3 (same values for both rows)
$values = @{
    "D" = 0.0399
    "E" = 0.0384
    "F" = 0.06
    "I" = 0
    "J" = 0
    "K" = 209
    "L" = 0.09745407068917281
    "M" = 47.03
    "N" = 0.02064802212758485
    "O" = 0.2250239234449761
    "P" = 43.6
    "Q" = 0.0191421170478992
    "R" = 0.2086124401913876
    "S" = 3.43
    "T" = 0.07293217095470975
    "U" = 865.8
    "V" = 0.3801202967906221
    "W" = 0.02992211659603711
    "X" = 0.224238770787501
    "Y" = -0.1943166541914639
    "Z" = 0.1140853911544722
    "AA" = 0
    "AB" = 0.0670005838861821
    "AC" = -0.0670005838861821
    "AD" = 15986.5
    "AE" = 0
    "AF" = 15986.5
    "AG" = 15120.7
    "AH" = 0.8752915539689666
    "AI" = 0.6454185036254704
    "AJ" = 0.8690856630494758
    "AK" = 0.6325752821774309
}

foreach ($row in 2..3) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
    # Cells AN and AP are removed entirely in the new version
    $ws.Range("AN$row").ClearContents()
    $ws.Range("AP$row").ClearContents()
}
